$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 594; this shifts rows 594:649 down to 595:650
# and Excel automatically extends the used range to A1:R650.
$ws.Rows.Item(594).Insert()

# Populate the newly inserted row 594 with its data.
$ws.Cells.Item(594, 1).Value = 5
$ws.Cells.Item(594, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(594, 3).Value = "Maule"
$ws.Cells.Item(594, 4).Value = 45166
$ws.Cells.Item(594, 5).Value = 7
$ws.Cells.Item(594, 6).Value = 100114014
$ws.Cells.Item(594, 7).Value = "Betarraga"
$ws.Cells.Item(594, 8).Value = "Sin especificar"
$ws.Cells.Item(594, 9).Value = "Primera"
$ws.Cells.Item(594, 10).Value = 5000
$ws.Cells.Item(594, 11).Value = 500
$ws.Cells.Item(594, 12).Value = 500
$ws.Cells.Item(594, 13).Value = 500
$ws.Cells.Item(594, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(594, 15).Value = "Región del Maule"
$ws.Cells.Item(594, 16).Value = 100
$ws.Cells.Item(594, 17).Value = 5
$ws.Cells.Item(594, 18).Value = "Hortaliza"
